$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new price rows before current row 107 (new data dated 2023-03-21)
$ws.Range("A107:A109").EntireRow.Insert()

    $ws.Range("A107").Value = 7
    $ws.Range("B107").Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Range("C107").Value = "Ñuble"
    $ws.Range("D107").Value = 45006
    $ws.Range("E107").Value = 16
    $ws.Range("F107").Value = 100112021
    $ws.Range("G107").Value = "Ají"
    $ws.Range("H107").Value = "Cacho cabra rojo"
    $ws.Range("I107").Value = "Primera"
    $ws.Range("J107").Value = 30
    $ws.Range("K107").Value = 15000
    $ws.Range("L107").Value = 15000
    $ws.Range("M107").Value = 15000
    $ws.Range("N107").Value = "$/saco 25 kilos"
    $ws.Range("O107").Value = "Región del Maule"
    $ws.Range("P107").Value = 600
    $ws.Range("Q107").Value = 25
    $ws.Range("R107").Value = "Hortaliza"

    $ws.Range("A108").Value = 7
    $ws.Range("B108").Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Range("C108").Value = "Ñuble"
    $ws.Range("D108").Value = 45006
    $ws.Range("E108").Value = 16
    $ws.Range("F108").Value = 100112021
    $ws.Range("G108").Value = "Ají"
    $ws.Range("H108").Value = "Cacho cabra rojo"
    $ws.Range("I108").Value = "Primera"
    $ws.Range("J108").Value = 25
    $ws.Range("K108").Value = 15000
    $ws.Range("L108").Value = 15000
    $ws.Range("M108").Value = 15000
    $ws.Range("N108").Value = "$/saco 25 kilos"
    $ws.Range("O108").Value = "Región del Maule"
    $ws.Range("P108").Value = 600
    $ws.Range("Q108").Value = 25
    $ws.Range("R108").Value = "Hortaliza"

    $ws.Range("A109").Value = 7
    $ws.Range("B109").Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Range("C109").Value = "Ñuble"
    $ws.Range("D109").Value = 45006
    $ws.Range("E109").Value = 16
    $ws.Range("F109").Value = 100112021
    $ws.Range("G109").Value = "Ají"
    $ws.Range("H109").Value = "Cristal"
    $ws.Range("I109").Value = "Primera"
    $ws.Range("J109").Value = 20
    $ws.Range("K109").Value = 15000
    $ws.Range("L109").Value = 15000
    $ws.Range("M109").Value = 15000
    $ws.Range("N109").Value = "$/saco 25 kilos"
    $ws.Range("O109").Value = "Región del Maule"
    $ws.Range("P109").Value = 600
    $ws.Range("Q109").Value = 25
    $ws.Range("R109").Value = "Hortaliza"

# Insert 3 more new price rows before current row 129 (new data dated 2023-03-20)
$ws.Range("A129:A131").EntireRow.Insert()

    $ws.Range("A129").Value = 7
    $ws.Range("B129").Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Range("C129").Value = "Ñuble"
    $ws.Range("D129").Value = 45005
    $ws.Range("E129").Value = 16
    $ws.Range("F129").Value = 100112021
    $ws.Range("G129").Value = "Ají"
    $ws.Range("H129").Value = "Cacho cabra rojo"
    $ws.Range("I129").Value = "Primera"
    $ws.Range("J129").Value = 40
    $ws.Range("K129").Value = 15000
    $ws.Range("L129").Value = 15000
    $ws.Range("M129").Value = 15000
    $ws.Range("N129").Value = "$/saco 25 kilos"
    $ws.Range("O129").Value = "Región del Maule"
    $ws.Range("P129").Value = 600
    $ws.Range("Q129").Value = 25
    $ws.Range("R129").Value = "Hortaliza"

    $ws.Range("A130").Value = 7
    $ws.Range("B130").Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Range("C130").Value = "Ñuble"
    $ws.Range("D130").Value = 45005
    $ws.Range("E130").Value = 16
    $ws.Range("F130").Value = 100112021
    $ws.Range("G130").Value = "Ají"
    $ws.Range("H130").Value = "Cacho cabra verde"
    $ws.Range("I130").Value = "Primera"
    $ws.Range("J130").Value = 30
    $ws.Range("K130").Value = 15000
    $ws.Range("L130").Value = 15000
    $ws.Range("M130").Value = 15000
    $ws.Range("N130").Value = "$/saco 25 kilos"
    $ws.Range("O130").Value = "Provincia de Diguillín"
    $ws.Range("P130").Value = 600
    $ws.Range("Q130").Value = 25
    $ws.Range("R130").Value = "Hortaliza"

    $ws.Range("A131").Value = 7
    $ws.Range("B131").Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Range("C131").Value = "Ñuble"
    $ws.Range("D131").Value = 45005
    $ws.Range("E131").Value = 16
    $ws.Range("F131").Value = 100112021
    $ws.Range("G131").Value = "Ají"
    $ws.Range("H131").Value = "Cristal"
    $ws.Range("I131").Value = "Primera"
    $ws.Range("J131").Value = 20
    $ws.Range("K131").Value = 15000
    $ws.Range("L131").Value = 15000
    $ws.Range("M131").Value = 15000
    $ws.Range("N131").Value = "$/saco 25 kilos"
    $ws.Range("O131").Value = "Región del Maule"
    $ws.Range("P131").Value = 600
    $ws.Range("Q131").Value = 25
    $ws.Range("R131").Value = "Hortaliza"

